$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 3461.6667
$ws.Range("I8").Value = 4054
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 12162
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -12023
$ws.Range("N8").Value = -1778
$ws.Range("H42").Value = 142.66667
$ws.Range("I42").Value = 39.833332
$ws.Range("J42").Value = 348.33334
$ws.Range("K42").Value = 119.499996
$ws.Range("L42").Value = 1045.00002
$ws.Range("M42").Value = 110.500004
$ws.Range("N42").Value = -1505.00002
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4909.1177
$ws.Range("I74").Value = 1502.6086
$ws.Range("J74").Value = 12031.818
$ws.Range("K74").Value = 1502.6086
$ws.Range("L74").Value = 12031.818
$ws.Range("M74").Value = -628.6086
$ws.Range("N74").Value = -13779.818
$ws.Range("H77").Value = 4909.1177
$ws.Range("I77").Value = 1502.6086
$ws.Range("J77").Value = 12031.818
$ws.Range("K77").Value = 7513.043
$ws.Range("L77").Value = 60159.09
$ws.Range("M77").Value = -3145.043
$ws.Range("N77").Value = -68895.09
$ws.Range("H102").Value = 1784.6666
$ws.Range("I102").Value = 1410
$ws.Range("J102").Value = 1972
$ws.Range("K102").Value = 1410
$ws.Range("L102").Value = 1972
$ws.Range("M102").Value = 212
$ws.Range("N102").Value = -5216
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 6666.6665
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 9500
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 9500
$ws.Range("M8").Value = -860
$ws.Range("N8").Value = -9780
$ws.Range("H10").Value = 1007
$ws.Range("I10").Value = 1007
$ws.Range("K10").Value = 1007
$ws.Range("M10").Value = -868
$ws.Range("H12").Value = 27711
$ws.Range("I12").Value = 26666.5
$ws.Range("J12").Value = 29800
$ws.Range("K12").Value = 26666.5
$ws.Range("L12").Value = 29800
$ws.Range("M12").Value = -26496.5
$ws.Range("N12").Value = -30140
$ws.Range("H13").Value = 1000000
$ws.Range("J13").Value = 1000000
$ws.Range("L13").Value = 1000000
$ws.Range("N13").Value = -1000278
$ws.Range("H14").Value = 38000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 38000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 38000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -38340
$ws.Range("H17").Value = 27500
$ws.Range("I17").Value = 27500
$ws.Range("K17").Value = 27500
$ws.Range("M17").Value = -27326
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H64").Value = 1931.8572
$ws.Range("J64").Value = 3666.3333
$ws.Range("L64").Value = 10998.9999
$ws.Range("N64").Value = -11538.9999
$ws.Range("H67").Value = 1931.8572
$ws.Range("J67").Value = 3666.3333
$ws.Range("L67").Value = 10998.9999
$ws.Range("N67").Value = -12870.9999
$ws.Range("H107").Value = 1684.7667
$ws.Range("I107").Value = 680
$ws.Range("J107").Value = 1817.4717
$ws.Range("K107").Value = 2040
$ws.Range("L107").Value = 5452.4151
$ws.Range("M107").Value = -120
$ws.Range("N107").Value = -9292.4151
$ws.Range("H131").Value = 3416.196
$ws.Range("I131").Value = 666.3333
$ws.Range("K131").Value = 1998.9999
$ws.Range("M131").Value = 3041.0001
$ws.Range("H133").Value = 11000
$ws.Range("I133").Value = 11000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 33000
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -27940
$ws.Range("N133").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 294428.56
$ws.Range("I3").Value = 1025000
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 1025000
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -1024884
$ws.Range("N3").Value = -2432
$ws.Range("H24").Value = 2516126.5
$ws.Range("I24").Value = 5025003
$ws.Range("K24").Value = 5025003
$ws.Range("M24").Value = -5024830
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H107").Value = 1428.8572
$ws.Range("I107").Value = 1428.8572
$ws.Range("K107").Value = 1428.8572
$ws.Range("M107").Value = 491.1428000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 39429
$ws.Range("J123").Value = 39429
$ws.Range("L123").Value = 39429
$ws.Range("N123").Value = -49229
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 22500
$ws.Range("J18").Value = 22500
$ws.Range("L18").Value = 22500
$ws.Range("N18").Value = -22846
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30480
$ws.Range("H22").Value = 22500
$ws.Range("J22").Value = 22500
$ws.Range("L22").Value = 22500
$ws.Range("N22").Value = -23086
$ws.Range("H52").Value = 37333.332
$ws.Range("I52").Value = 100000
$ws.Range("J52").Value = 6000
$ws.Range("K52").Value = 100000
$ws.Range("L52").Value = 6000
$ws.Range("M52").Value = -99774
$ws.Range("N52").Value = -6452
$ws.Range("H62").Value = 16702814
$ws.Range("J62").Value = 29221.75
$ws.Range("L62").Value = 29221.75
$ws.Range("N62").Value = -30469.75
$ws.Range("H65").Value = 16702814
$ws.Range("J65").Value = 29221.75
$ws.Range("L65").Value = 146108.75
$ws.Range("N65").Value = -152348.75
$ws.Range("H75").Value = 143529.5
$ws.Range("I75").Value = 257559
$ws.Range("J75").Value = 29500
$ws.Range("K75").Value = 257559
$ws.Range("L75").Value = 29500
$ws.Range("M75").Value = -256623
$ws.Range("N75").Value = -31372
$ws.Range("H78").Value = 143529.5
$ws.Range("I78").Value = 257559
$ws.Range("J78").Value = 29500
$ws.Range("K78").Value = 772677
$ws.Range("L78").Value = 88500
$ws.Range("M78").Value = -767997
$ws.Range("N78").Value = -97860
$ws.Range("H81").Value = 1616722.8
$ws.Range("I81").Value = 1751282.9
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 3502565.8
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -3501504.8
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1616722.8
$ws.Range("I84").Value = 1751282.9
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 17512829
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -17507525
$ws.Range("N84").Value = -30608
$ws.Range("H96").Value = 55556508
$ws.Range("I96").Value = 90910290
$ws.Range("J96").Value = 571.4286
$ws.Range("K96").Value = 90910290
$ws.Range("L96").Value = 571.4286
$ws.Range("M96").Value = -90908917
$ws.Range("N96").Value = -3317.4286
$ws.Range("H132").Value = 17859788
$ws.Range("I132").Value = 25002154
$ws.Range("K132").Value = 75006462
$ws.Range("M132").Value = -75003932
$ws.Range("H136").Value = 8360340
$ws.Range("I136").Value = 23881604
$ws.Range("K136").Value = 71644812
$ws.Range("M136").Value = -71642262
